$wb = $excel.ActiveWorkbook

$wsChart = $wb.Worksheets.Item("Chart")
$wsTable = $wb.Worksheets.Item("Table")

# Update the status string "Not Started" -> "Started" on the Table sheet
$wsTable.Range("B2").Value = "Started"

# Update D45 on the Chart sheet from blank text to numeric 0
$wsChart.Range("D45").Value = 0
